$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "59.496.88"
$ws.Range("E2").Value = "  +2.94%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.984.97"
$ws.Range("E3").Value = "  +1.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "566.34"
$ws.Range("E5").Value = "  +2.42%  "

# Row 6 - Solana
Set-TextValue "D6" "138.62"
$ws.Range("E6").Value = "  +3.84%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 - XRP
Set-TextValue "D8" "0.520"
$ws.Range("E8").Value = "  +1.19%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.975.81"
$ws.Range("E9").Value = "  +1.30%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +3.69%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.38"
$ws.Range("E11").Value = "  +12.38%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.59%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000230"
$ws.Range("E13").Value = "  +4.38%  "

# Row 14 - Avalanche
Set-TextValue "D14" "33.77"
$ws.Range("E14").Value = "  +3.01%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.21%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.477.69"
$ws.Range("E16").Value = "  +1.30%  "

# Row 17 - Polkadot
Set-TextValue "D17" "7.06"
$ws.Range("E17").Value = "  +1.88%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.981.04"
$ws.Range("E18").Value = "  +1.23%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "59.485.73"
$ws.Range("E19").Value = "  +2.96%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "437.05"
$ws.Range("E20").Value = "  +5.17%  "

# Row 21 - Chainlink
Set-TextValue "D21" "13.64"
$ws.Range("E21").Value = "  +2.18%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.723"
$ws.Range("E22").Value = "  +3.67%  "

# Row 23 - was Uniswap, now InternetComputer(DFINITY)
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D23" "13.33"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24 - was InternetComputer(DFINITY), now Uniswap
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D24" "7.03"
$ws.Range("E24").Value = "  +0.45%  "

# Row 25 - Litecoin
Set-TextValue "D25" "80.04"
$ws.Range("E25").Value = "  +1.32%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.08%  "

# Row 27 - ImmutableX
Set-TextValue "D27" "2.23"
$ws.Range("E27").Value = "  +10.16%  "

# Row 28 - FirstDigitalUSD
$ws.Range("E28").Value = "  -0.03%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +2.51%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.76"
$ws.Range("E30").Value = "  +3.39%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "25.77"
$ws.Range("E31").Value = "  +1.32%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "6.22"
$ws.Range("E32").Value = "  +4.48%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.105"
$ws.Range("E33").Value = "  +9.31%  "

# Row 34 - PEPE
Set-TextValue "D34" "0.0₃0777"
$ws.Range("E34").Value = "  +11.39%  "

# Row 35 - Filecoin
Set-TextValue "D35" "5.92"
$ws.Range("E35").Value = "  +4.24%  "

# Row 36 - Mantle
Set-TextValue "D36" "0.980"
$ws.Range("E36").Value = "  +3.01%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  +0.53%  "

# Row 38 - OKB
Set-TextValue "D38" "48.67"
$ws.Range("E38").Value = "  +0.66%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  +2.82%  "

# Row 41 - Bittensor
Set-TextValue "D41" "399.05"
$ws.Range("E41").Value = "  +4.67%  "

# Row 42 - VeChain
Set-TextValue "D42" "0.0352"
$ws.Range("E42").Value = "  +1.40%  "

# Row 43 - Maker
Set-TextValue "D43" "2.731.47"
$ws.Range("E43").Value = "  +0.93%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  -2.08%  "

# Row 45 - TheGraph
Set-TextValue "D45" "0.251"
$ws.Range("E45").Value = "  +6.02%  "

# Row 46 - Arweave
Set-TextValue "D46" "35.17"
$ws.Range("E46").Value = "  +21.48%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  -0.04%  "

# Row 48 - Monero
Set-TextValue "D48" "122.11"
$ws.Range("E48").Value = "  -1.84%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  +2.23%  "

# Row 50 - Fetch.AI
$ws.Range("E50").Value = "  +1.93%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "23.30"
$ws.Range("E51").Value = "  +1.83%  "
